$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels for B1/C1 ("long" <-> "lat")
$ws.Range("B1").Value = "lat"
$ws.Range("C1").Value = "long"

# Uppercase the province name in column G for rows 2-4
$ws.Range("G2").Value = "JAWA BARAT"
$ws.Range("G3").Value = "JAWA BARAT"
$ws.Range("G4").Value = "JAWA BARAT"

# Update coordinates for row 3
$ws.Range("B3").Value = -6.5521228000000002
$ws.Range("C3").Value = 106.7535814

# Update the selected cell
$ws.Range("C4").Select()
